# Rename arguments of measure add_temperature_setpoints to conform to OS
# style and to improve clarity.
#
# This adds a new "add_temperature_setpoints" measure block (13 argument
# rows) to the bottom of the arguments table, and fixes the "Required"
# flag of the pre-existing zone_cooling_temp_sched_holiday argument
# (row 22, column G, measure add_detailed_hvac) from True to False.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New measure block: add_temperature_setpoints (rows 30-42)
# Columns: A=Measure, B=Argument name, C=Argument display name,
#          D=Type, E=Unit, F=Description, G=Required (True/False)
$measure = "add_temperature_setpoints"

$rows = @(
    @{ B = "zone_heating_temp_sched_weekday"; C = "Heating temperature schedule for weekdays"; D = "String"; F = "Schedule for the heating temperature setpoints for weekdays.";  G = "True"  },
    @{ B = "zone_heating_temp_sched_saturday"; C = "Heating temperature schedule for saturday";  D = "String"; F = "Schedule for the heating temperature setpoints for saturday.";  G = "True"  },
    @{ B = "zone_heating_temp_sched_sunday";   C = "Heating temperature schedule for sunday";    D = "String"; F = "Schedule for the heating temperature setpoints for sunday.";    G = "True"  },
    @{ B = "zone_heating_temp_sched_holiday";  C = "Heating temperature schedule for holidays";  D = "String"; F = "Schedule for the heating temperature setpoints for holidays.";  G = "False" },
    @{ B = "zone_cooling_temp_sched_weekday";  C = "Cooling temperature schedule for weekdays";  D = "String"; F = "Schedule for the cooling temperature setpoints for weekdays.";  G = "True"  },
    @{ B = "zone_cooling_temp_sched_saturday"; C = "Cooling temperature schedule for saturday";  D = "String"; F = "Schedule for the cooling temperature setpoints for saturday.";  G = "True"  },
    @{ B = "zone_cooling_temp_sched_sunday";   C = "Cooling temperature schedule for sunday";    D = "String"; F = "Schedule for the cooling temperature setpoints for sunday.";    G = "True"  },
    @{ B = "zone_cooling_temp_sched_holiday";  C = "Cooling temperature schedule for holiday";   D = "String"; F = "Schedule for the cooling temperature setpoints for holidays.";  G = "False" },
    @{ B = "holidays";                         C = "Holiday definition";                          D = "String"; F = "Definition of holidays";                                       G = "False" },
    @{ B = "heating_temp_selection";           C = "Selection of heating temperature";            D = "String"; F = "(Export only) Selection of heating temperature";                G = "False" },
    @{ B = "cooling_temp_selection";           C = "Selection of cooling temperature";            D = "String"; F = "(Export only) Selection of cooling temperature";                G = "False" },
    @{ B = "is_custom_heating";                C = "Is custom heating temperature";               D = "Bool";   F = "(Export only) Flag whether the selected heating temperature is a custom schedule"; G = "False" },
    @{ B = "is_custom_cooling";                C = "Is custom cooling temperature";               D = "Bool";   F = "(Export only) Flag whether the selectedcooling temperature is a custom schedule";  G = "False" }
)

# Cells that already contain the literal text "True"/"False" stored as
# shared strings (not as native booleans). Assigning the strings "True"/
# "False" directly via .Value would get auto-converted to a boolean cell
# type, so instead we copy these known-good text cells wherever a
# True/False label is needed.
$trueCell  = $ws.Range("G2")
$falseCell = $ws.Range("G8")

$r = 30
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $measure
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 6).Value = $row.F

    if ($row.G -eq "True") {
        $trueCell.Copy($ws.Cells.Item($r, 7)) | Out-Null
    } else {
        $falseCell.Copy($ws.Cells.Item($r, 7)) | Out-Null
    }

    $r = $r + 1
}

# --- Fix existing row: add_detailed_hvac -> zone_cooling_temp_sched_holiday
# "Required" column (G) should be False, not True.
$falseCell.Copy($ws.Range("G22")) | Out-Null

# --- Resize columns to fit the new content, like Excel's AutoFit.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null
$ws.Columns.Item(6).AutoFit() | Out-Null
$ws.Columns.Item(8).AutoFit() | Out-Null

# --- Move/adjust the active selection to A43, just past the new data,
# matching the end-state selection recorded in the workbook.
$ws.Range("A43").Select() | Out-Null
